# Compare Plan Test Scenarios
# Adds new Compare-Plan related test rows to the "Platform" sheet and
# flips the existing VerifyComparePlanMenuLauncherTest status to "N".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Platform")

# Existing row 11 (VerifyComparePlanMenuLauncherTest) now reports "N"
$ws.Range("D11").Value = "N"

# New test-case rows, following the same Android/Oneplus pattern used
# throughout the sheet.
$newRows = @(
    @{ Row = 12; TC = "ComparePlanCloseDrawerTest" },
    @{ Row = 13; TC = "CloseComparePlanTest" },
    @{ Row = 14; TC = "ContinueComparePlanTest" }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Range("A$row").Value = $r.TC
    $ws.Range("B$row").Value = "Android"
    $ws.Range("C$row").Value = "Oneplus"
    $ws.Range("C$row").HorizontalAlignment = -4108
    $ws.Range("C$row").VerticalAlignment = -4108
    $ws.Range("D$row").Value = "Y"
}

# Match the selection left behind by the edit (active cell D12).
$ws.Range("D12").Select() | Out-Null
